$d = $word.ActiveDocument

# Find the paragraph "Marketing Strategy and Data-Driven Insights" under the
# Siege Analytics / Partner role, and insert three new bullet paragraphs
# right after it (before the "Conducted comprehensive..." bullet).

$targetText = "Marketing Strategy and Data-Driven Insights"
$found = $false

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $targetText) {
        $found = $true

        # Position right before the paragraph mark that ends this paragraph.
        $insertPos = $p.Range.End - 1
        $r = $d.Range($insertPos, $insertPos)

        $newText = "`r" + "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters" `
                 + "`r" + "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States" `
                 + "`r" + "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"

        $r.InsertAfter($newText)
        break
    }
}

if (-not $found) {
    throw "Could not find the 'Marketing Strategy and Data-Driven Insights' paragraph"
}
